$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update July (row 8) 2021 value
$ws.Range("H8").Value = 153

# Update August row label (row 9) to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-17)"

# Update August (row 9) year values
$ws.Range("B9").Value = 20
$ws.Range("D9").Value = 43
$ws.Range("E9").Value = 27
$ws.Range("F9").Value = 24
$ws.Range("G9").Value = 111
$ws.Range("H9").Value = 93

# Update Total (row 10) year values
$ws.Range("B10").Value = 182
$ws.Range("D10").Value = 508
$ws.Range("E10").Value = 452
$ws.Range("F10").Value = 328
$ws.Range("G10").Value = 732
